# "run on more subjects"
#
# The "In Lab" sheet had a small participant-info block typed into columns
# K:N (columns H:J were left blank as a visual gutter between the main
# timing table in A:G and that block). To make room to log more subjects,
# two of those blank spacer columns (H:I) were deleted, which shifts the
# K:N block left by two columns to I:L - pulling it closer to the main
# table so there's more room to the right for additional subject columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("In Lab")

# Delete the two blank columns that separated the timing table (A:G) from
# the participant-info block (K:N). Excel automatically shifts everything
# to the right of the deletion (K:N) left by two columns, becoming I:L,
# and updates column widths / cell references / dimension accordingly.
$ws.Range("H1:I1").EntireColumn.Delete()

# Reflect the selection that was left after the edit: the whole of the
# (now relocated) first info column, I, with the active cell at the top.
$ws.Activate()
$excel.Goto($ws.Columns("I:I"), $true)
